# Regenerate save_data to use K (strikeouts) instead of Strike# in column G,
# writing the newly calculated s_vals for each saved game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(0,2,3,0,1,2,0,1,1,1,0,0,1,1,1,0,0,2,2,3,1,1,1,1,1,1,1,2,1,3,0,1,0,1,1,1,0,1,3,1,0,1,1,0,0,2,1,0,0,1,1,1,1,1,2,2,2,2,0,0,1,1,2,1,1,3,1,1,1,1,2,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
